$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric values for rows 2-6 (columns D..AJ)
# Row 2
$ws.Range("D2").Value = 30244
$ws.Range("E2").Value = -1601
$ws.Range("F2").Value = -1601
$ws.Range("G2").Value = -14167
$ws.Range("H2").Value = -12812
$ws.Range("I2").Value = -12746
$ws.Range("J2").Value = -66
$ws.Range("K2").Value = 33939
$ws.Range("L2").Value = 33244
$ws.Range("M2").Value = 695
$ws.Range("N2").Value = 515
$ws.Range("O2").Value = 180
$ws.Range("P2").Value = 3711
$ws.Range("Q2").Value = -318
$ws.Range("R2").Value = 164
$ws.Range("S2").Value = -104
$ws.Range("T2").Value = 673
$ws.Range("U2").Value = -991
$ws.Range("V2").Value = 22735
$ws.Range("W2").Value = -5.29
$ws.Range("X2").Value = -42.36
$ws.Range("Y2").Value = -199.87
$ws.Range("Z2").Value = -29.67
$ws.Range("AA2").Value = 4783.09
$ws.Range("AB2").Value = -82.38
$ws.Range("AC2").Value = -5631863
$ws.Range("AD2").Value = -0.01
$ws.Range("AE2").Value = 223294
$ws.Range("AF2").Value = 0.3
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 182752

# Row 3
$ws.Range("D3").Value = 23208
$ws.Range("E3").Value = 786
$ws.Range("F3").Value = 786
$ws.Range("G3").Value = -488
$ws.Range("H3").Value = -481
$ws.Range("I3").Value = -462
$ws.Range("J3").Value = -19
$ws.Range("K3").Value = 28905
$ws.Range("L3").Value = 28871
$ws.Range("M3").Value = 34
$ws.Range("N3").Value = 16
$ws.Range("O3").Value = 18
$ws.Range("P3").Value = 1582
$ws.Range("Q3").Value = 903
$ws.Range("R3").Value = 2024
$ws.Range("S3").Value = -2570
$ws.Range("T3").Value = 341
$ws.Range("U3").Value = 562
$ws.Range("V3").Value = 21809
$ws.Range("W3").Value = 3.39
$ws.Range("X3").Value = -2.07
$ws.Range("Y3").Value = -174.04
$ws.Range("Z3").Value = -1.53
$ws.Range("AA3").Value = 84510.24
$ws.Range("AB3").Value = -104.39
$ws.Range("AC3").Value = -110540
$ws.Range("AD3").Value = -0.57
$ws.Range("AE3").Value = 3673
$ws.Range("AF3").Value = 17.13
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 388824

# Row 4
$ws.Range("D4").Value = 23280
$ws.Range("E4").Value = 1436
$ws.Range("F4").Value = 1436
$ws.Range("G4").Value = -638
$ws.Range("H4").Value = -727
$ws.Range("I4").Value = -727
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 28232
$ws.Range("L4").Value = 26970
$ws.Range("M4").Value = 1263
$ws.Range("N4").Value = 1244
$ws.Range("O4").Value = 18
$ws.Range("P4").Value = 1808
$ws.Range("Q4").Value = 3175
$ws.Range("R4").Value = -289
$ws.Range("S4").Value = -2163
$ws.Range("T4").Value = 290
$ws.Range("U4").Value = 2884
$ws.Range("V4").Value = 19935
$ws.Range("W4").Value = 6.17
$ws.Range("X4").Value = -3.12
$ws.Range("Y4").Value = -115.44
$ws.Range("Z4").Value = -2.54
$ws.Range("AA4").Value = 2136.01
$ws.Range("AB4").Value = -34.11
$ws.Range("AC4").Value = -50193
$ws.Range("AD4").Value = -0.87
$ws.Range("AE4").Value = 62313
$ws.Range("AF4").Value = 0.7
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 1944082

# Row 5
$ws.Range("D5").Value = 25962
$ws.Range("E5").Value = -118
$ws.Range("F5").Value = -118
$ws.Range("G5").Value = -1636
$ws.Range("H5").Value = -1607
$ws.Range("I5").Value = -1608
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 26544
$ws.Range("L5").Value = 25047
$ws.Range("M5").Value = 1497
$ws.Range("N5").Value = 1478
$ws.Range("O5").Value = 19
$ws.Range("P5").Value = 1808
$ws.Range("Q5").Value = 1831
$ws.Range("R5").Value = 97
$ws.Range("S5").Value = -1988
$ws.Range("T5").Value = 280
$ws.Range("U5").Value = 1551
$ws.Range("V5").Value = 17553
$ws.Range("W5").Value = -0.45
$ws.Range("X5").Value = -6.19
$ws.Range("Y5").Value = -118.17
$ws.Range("Z5").Value = -5.87
$ws.Range("AA5").Value = 1673.42
$ws.Range("AB5").Value = -18.8
$ws.Range("AC5").Value = -80525
$ws.Range("AD5").Value = -0.28
$ws.Range("AE5").Value = 74029
$ws.Range("AF5").Value = 0.3
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 1944082

# Row 6
$ws.Range("D6").Value = 25451
$ws.Range("E6").Value = -656
$ws.Range("F6").Value = -656
$ws.Range("G6").Value = -1177
$ws.Range("H6").Value = -1183
$ws.Range("I6").Value = -1179
$ws.Range("K6").Value = 25689
$ws.Range("L6").Value = 25080
$ws.Range("M6").Value = 609
$ws.Range("N6").Value = 605
$ws.Range("P6").Value = 1920
$ws.Range("Q6").Value = -82
$ws.Range("R6").Value = -327
$ws.Range("S6").Value = -108
$ws.Range("T6").Value = 299
$ws.Range("U6").Value = -382
$ws.Range("V6").Value = 17417
$ws.Range("W6").Value = -2.58
$ws.Range("X6").Value = -4.65
$ws.Range("Y6").Value = -113.18
$ws.Range("Z6").Value = -4.53
$ws.Range("AA6").Value = 4120.01
$ws.Range("AB6").Value = -70.18
$ws.Range("AC6").Value = -34755
$ws.Range("AD6").Value = -0.55
$ws.Range("AE6").Value = 15675
$ws.Range("AF6").Value = 1.23
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 3808897

# Rows 7-9: clear all data columns (D:AJ), keeping only A, B, C
$clearRange = $ws.Range("D7:AJ7,D8:AJ8,D9:AJ9")
foreach ($area in $clearRange.Areas) {
    $area.ClearContents()
}
